$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, "01/04/2015", 889, 0.45),
    @(3, "01/04/2015", 2327, -14.86),
    @(4, "01/04/2015", 1857, -20.51),
    @(5, "01/04/2015", 1857, -9.369999999999999),
    @(6, "01/04/2015", 2174, 6.52),
    @(7, "01/04/2015", 1450, -15.65),
    @(8, "01/04/2015", 3109, 3.12),
    @(9, "01/04/2015", 4023, -13.78),
    @(10, "01/04/2016", 814, -8.44),
    @(11, "01/04/2016", 2441, 4.9),
    @(12, "01/04/2016", 1665, -10.34),
    @(13, "01/04/2016", 1910, 2.85),
    @(14, "01/04/2016", 2228, 2.48),
    @(15, "01/04/2016", 1565, 7.93),
    @(16, "01/04/2016", 2947, -5.21),
    @(17, "01/04/2016", 4558, 13.3),
    @(18, "01/04/2017", 908, 11.55),
    @(19, "01/04/2017", 2177, -10.82),
    @(20, "01/04/2017", 1812, 8.83),
    @(21, "01/04/2017", 1779, -6.86),
    @(22, "01/04/2017", 1869, -16.11),
    @(23, "01/04/2017", 1600, 2.24),
    @(24, "01/04/2017", 2833, -3.87),
    @(25, "01/04/2017", 4460, -2.15),
    @(26, "01/04/2018", 990, 9.029999999999999),
    @(27, "01/04/2018", 1820, -16.4),
    @(28, "01/04/2018", 1645, -9.220000000000001),
    @(29, "01/04/2018", 1956, 9.949999999999999),
    @(30, "01/04/2018", 1549, -17.12),
    @(31, "01/04/2018", 1491, -6.81),
    @(32, "01/04/2018", 2704, -4.55),
    @(33, "01/04/2018", 3865, -13.34),
    @(34, "01/04/2019", 759, -23.33),
    @(35, "01/04/2019", 1660, -8.789999999999999),
    @(36, "01/04/2019", 1591, -3.28),
    @(37, "01/04/2019", 1716, -12.27),
    @(38, "01/04/2019", 1629, 5.16),
    @(39, "01/04/2019", 1577, 5.77),
    @(40, "01/04/2019", 3041, 12.46),
    @(41, "01/04/2019", 4002, 3.54),
    @(42, "01/04/2020", 802, 5.67),
    @(43, "01/04/2020", 2004, 20.72),
    @(44, "01/04/2020", 1574, -1.07),
    @(45, "01/04/2020", 1990, 15.97),
    @(46, "01/04/2020", 1744, 7.06),
    @(47, "01/04/2020", 1787, 13.32),
    @(48, "01/04/2020", 2989, -1.71),
    @(49, "01/04/2020", 4004, 0.05),
    @(50, "01/04/2021", 809, 0.87),
    @(51, "01/04/2021", 2740, 36.73),
    @(52, "01/04/2021", 1537, -2.35),
    @(53, "01/04/2021", 1560, -21.61),
    @(54, "01/04/2021", 1707, -2.12),
    @(55, "01/04/2021", 1280, -28.37),
    @(56, "01/04/2021", 3353, 12.18),
    @(57, "01/04/2021", 4241, 5.92),
    @(58, "01/04/2022", 831, 2.72),
    @(59, "01/04/2022", 1941, -29.16),
    @(60, "01/04/2022", 1258, -18.15),
    @(61, "01/04/2022", 1554, -0.38),
    @(62, "01/04/2022", 1478, -13.42),
    @(63, "01/04/2022", 1192, -6.87),
    @(64, "01/04/2022", 2734, -18.46),
    @(65, "01/04/2022", 3691, -12.97),
    @(66, "01/04/2023", 901, 8.42),
    @(67, "01/04/2023", 2292, 18.08),
    @(68, "01/04/2023", 1416, 12.56),
    @(69, "01/04/2023", 1666, 7.21),
    @(70, "01/04/2023", 1938, 31.12),
    @(71, "01/04/2023", 1163, -2.43),
    @(72, "01/04/2023", 2784, 1.83),
    @(73, "01/04/2023", 3785, 2.55),
    @(74, "01/04/2024", 893, -0.89),
    @(75, "01/04/2024", 1872, -18.32),
    @(76, "01/04/2024", 1588, 12.15),
    @(77, "01/04/2024", 2092, 25.57),
    @(78, "01/04/2024", 2354, 21.47),
    @(79, "01/04/2024", 1322, 13.67),
    @(80, "01/04/2024", 2762, -0.79),
    @(81, "01/04/2024", 3737, -1.27)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $cVal = $entry[1]
    $dVal = $entry[2]
    $eVal = $entry[3]

    $cCell = $ws.Cells.Item($r, 3)
    $cCell.NumberFormat = "@"
    $cCell.Value = $cVal
    $cCell.Style = "Normal"

    $ws.Cells.Item($r, 4).Value = $dVal
    $ws.Cells.Item($r, 5).Value = $eVal
}